$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.742600083351135
$ws.Range("B1").Value = 3.500783681869507
$ws.Range("C1").Value = 2.934829711914062
$ws.Range("D1").Value = 2.018598556518555
$ws.Range("E1").Value = 1.172240734100342
